# WelcomeIntros.pptx update — gitbook 2025-06-05 15:01:49
#
# Changes:
#  1. Agenda slide ("Agenda") bullet list gets a new "Photos" bullet
#     right after "Introductions".
#  2. Two new title slides are inserted right before the final
#     "THANK YOU" slide:
#       - "Introductions"
#       - "Photos"

$p = $ppt.ActivePresentation

# --- 1. Agenda slide: add "Photos" bullet under "Introductions" -----------
$agenda = $p.Slides.Item(2)
$contentPlaceholder = $agenda.Shapes.Item(2)
$contentPlaceholder.TextFrame.TextRange.Text = "Introductions`rPhotos`r"

# --- 2. Insert new section-title slides before the closing slide ----------
# Current order: 1 PEEKE Project 2025 | 2 Agenda | 3 THANK YOU
# Target order:  1 PEEKE Project 2025 | 2 Agenda | 3 Introductions | 4 Photos | 5 THANK YOU
$introSlide = $p.Slides.Add(3, 1)
$introSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Introductions"

$photosSlide = $p.Slides.Add(4, 1)
$photosSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Photos"
